$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J (year 2021) values, mirroring column I (year 2020) rows 4-14.
$values = @{
    4  = 2021
    5  = 24.4
    6  = 45.7
    7  = 38
    8  = 51.3
    9  = 51.5
    10 = 13
    11 = 36.4
    12 = 27
    13 = 2.7
    14 = 40.4
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Cells.Item($row, 9)   # column I
    $dstCell = $ws.Cells.Item($row, 10)  # column J

    # Copy formatting (number format, font, borders, etc.) from column I so the
    # new column J cell is styled the same way as the existing data column.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats

    $dstCell.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Row 3 height changes from 18 to 13.5.
$ws.Rows.Item(3).RowHeight = 13.5

# Update the active selection to K18 as recorded after the edit.
$ws.Range("K18").Select()
